$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 4).Value = "23.189.62"
$ws.Cells.Item(2, 5).Value = "  -3.18%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.607.62"
$ws.Cells.Item(3, 5).Value = "  -2.72%  "

# Row 4
Set-TextValue 4 4 "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.07%  "

# Row 5
$ws.Cells.Item(5, 5).Value = "  +0.02%  "

# Row 6
Set-TextValue 6 4 "302.91"
$ws.Cells.Item(6, 5).Value = "  -2.12%  "

# Row 7
Set-TextValue 7 4 "0.3766"
$ws.Cells.Item(7, 5).Value = "  -3.36%  "

# Row 8
Set-TextValue 8 4 "0.3649"
$ws.Cells.Item(8, 5).Value = "  -4.73%  "

# Row 9
Set-TextValue 9 4 "48.73"
$ws.Cells.Item(9, 5).Value = "  -4.98%  "

# Row 10
Set-TextValue 10 4 "1.001"
$ws.Cells.Item(10, 5).Value = "  +0.00%  "

# Row 11
Set-TextValue 11 4 "1.269"
$ws.Cells.Item(11, 5).Value = "  -6.18%  "

# Row 12
Set-TextValue 12 4 "0.08077"
$ws.Cells.Item(12, 5).Value = "  -4.24%  "

# Row 13
Set-TextValue 13 4 "23.00"
$ws.Cells.Item(13, 5).Value = "  -3.77%  "

# Row 14
Set-TextValue 14 4 "6.572"
$ws.Cells.Item(14, 5).Value = "  -7.44%  "

# Row 15
Set-TextValue 15 4 "7.605"
$ws.Cells.Item(15, 5).Value = "  -3.61%  "

# Row 16
Set-TextValue 16 4 "0.00001267"
$ws.Cells.Item(16, 5).Value = "  -3.55%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "1.609.98"
$ws.Cells.Item(17, 5).Value = "  -2.68%  "

# Row 18
Set-TextValue 18 4 "91.51"
$ws.Cells.Item(18, 5).Value = "  -3.21%  "

# Row 19
Set-TextValue 19 4 "0.06784"
$ws.Cells.Item(19, 5).Value = "  -2.94%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -6.95%  "

# Row 21
Set-TextValue 21 4 "6.578"
$ws.Cells.Item(21, 5).Value = "  -5.08%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +0.04%  "

# Row 23
Set-TextValue 23 4 "13.10"
$ws.Cells.Item(23, 5).Value = "  -4.47%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "23.201.66"
$ws.Cells.Item(24, 5).Value = "  -3.16%  "

# Row 25
Set-TextValue 25 4 "2.357"
$ws.Cells.Item(25, 5).Value = "  -4.30%  "

# Row 26
Set-TextValue 26 4 "2.921"
$ws.Cells.Item(26, 5).Value = "  -2.14%  "

# Row 27
Set-TextValue 27 4 "21.13"
$ws.Cells.Item(27, 5).Value = "  -4.41%  "

# Row 28
Set-TextValue 28 4 "150.34"
$ws.Cells.Item(28, 5).Value = "  -0.55%  "

# Row 29
Set-TextValue 29 4 "5.262"
$ws.Cells.Item(29, 5).Value = "  -3.36%  "

# Row 30
Set-TextValue 30 4 "132.68"
$ws.Cells.Item(30, 5).Value = "  -4.62%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -2.92%  "

# Row 32
Set-TextValue 32 4 "6.944"
$ws.Cells.Item(32, 5).Value = "  -11.19%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "1.786.63"

# Row 34
Set-TextValue 34 4 "0.9784"
$ws.Cells.Item(34, 5).Value = "  -6.21%  "

# Row 35
Set-TextValue 35 4 "0.07727"
$ws.Cells.Item(35, 5).Value = "  -4.37%  "

# Row 36
Set-TextValue 36 4 "0.02783"
$ws.Cells.Item(36, 5).Value = "  -5.86%  "

# Row 37
Set-TextValue 37 4 "6.275"
$ws.Cells.Item(37, 5).Value = "  -6.85%  "

# Row 38
Set-TextValue 38 4 "0.2554"
$ws.Cells.Item(38, 5).Value = "  -4.86%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "Stellar"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue 39 4 "0.08856"
$ws.Cells.Item(39, 5).Value = "  -3.08%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "FraxShare"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue 40 4 "10.08"
$ws.Cells.Item(40, 5).Value = "  -7.30%  "

# Row 41
Set-TextValue 41 4 "1.394"
$ws.Cells.Item(41, 5).Value = "  -2.27%  "

# Row 42
Set-TextValue 42 4 "0.7158"
$ws.Cells.Item(42, 5).Value = "  -5.24%  "

# Row 43
Set-TextValue 43 4 "12.79"
$ws.Cells.Item(43, 5).Value = "  -4.96%  "

# Row 44
Set-TextValue 44 4 "15.86"
$ws.Cells.Item(44, 5).Value = "  -2.98%  "

# Row 45
Set-TextValue 45 4 "0.6601"
$ws.Cells.Item(45, 5).Value = "  -4.96%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "Frax"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue 46 4 "1.000"
$ws.Cells.Item(46, 5).Value = "  +0.03%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "NEARProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 47 4 "2.300"
$ws.Cells.Item(47, 5).Value = "  -6.29%  "

# Row 48
Set-TextValue 48 4 "3.981"
$ws.Cells.Item(48, 5).Value = "  -2.68%  "

# Row 49
Set-TextValue 49 4 "0.08014"
$ws.Cells.Item(49, 5).Value = "  -3.29%  "

# Row 50
Set-TextValue 50 4 "131.30"
$ws.Cells.Item(50, 5).Value = "  -2.58%  "

# Row 51
Set-TextValue 51 4 "1.168"
$ws.Cells.Item(51, 5).Value = "  -3.34%  "
